$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.515.64"
$ws.Range("E2").Value = "  -6.84%  "

$ws.Range("D3").Value = "2.585.00"
$ws.Range("E3").Value = "  -1.37%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "301.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.54%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "96.29"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.08%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.580"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.28%  "

$ws.Range("E8").Value = "  +0.25%  "

$ws.Range("E9").Value = "  -4.36%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.60"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -8.29%  "

$ws.Range("E11").Value = "  -4.19%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.75"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -6.02%  "

$ws.Range("D13").Value = "2.984.91"
$ws.Range("E13").Value = "  -0.94%  "

$ws.Range("E14").Value = "  +0.88%  "

$ws.Range("D15").Value = "2.586.57"
$ws.Range("E15").Value = "  -1.22%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.889"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.86%  "

$ws.Range("E17").Value = "  -5.28%  "

$ws.Range("D18").Value = "43.517.87"

$ws.Range("D19").Value = "0.0₃0981"
$ws.Range("E19").Value = "  -4.32%  "

$ws.Range("E20").Value = "  -1.76%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.01%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.84"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.98%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "265.17"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.03%  "

$ws.Range("E24").Value = "  -4.28%  "

$ws.Range("E25").Value = "  +1.01%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "29.22"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.88%  "

$ws.Range("E27").Value = "  +0.23%  "

$ws.Range("E28").Value = "  -4.17%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "37.90"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.81%  "

$ws.Range("E30").Value = "  -6.69%  "

$ws.Range("E31").Value = "  -6.26%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.62"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.28%  "

$ws.Range("E33").Value = "  -1.97%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "152.27"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.19%  "

$ws.Range("E35").Value = "  -2.51%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0811"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.62%  "

$ws.Range("E37").Value = "  -4.60%  "

$ws.Range("E38").Value = "  -2.51%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "24.21"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.91%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "16.63"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.57%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.60"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.10%  "

$ws.Range("E42").Value = "  -5.54%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.88"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.88%  "

$ws.Range("D44").Value = "2.041.48"
$ws.Range("E44").Value = "  -3.99%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.997"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.13%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "87.77"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.41%  "

$ws.Range("E47").Value = "  -5.06%  "

$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.62"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.16%  "

$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "2.840.63"
$ws.Range("E49").Value = "  -1.12%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "105.60"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.08%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.190"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.92%  "
